$d = $word.ActiveDocument

$replacements = @(
    @("57×31=", "54×32="),
    @("55×80=", "33×77="),
    @("16×49=", "59×80="),
    @("58×32=", "48×24="),
    @("48×76=", "83×64="),
    @("73×18=", "70×99="),
    @("11×31=", "89×44="),
    @("83×40=", "27×83="),
    @("63×73=", "15×71="),
    @("69×40=", "16×88="),
    @("66×73=", "70×72="),
    @("18×40=", "81×96="),
    @("77×45=", "99×78="),
    @("63×42=", "81×87="),
    @("40×46=", "93×55="),
    @("19×84=", "79×53="),
    @("39×64=", "97×93="),
    @("46×22=", "28×20="),
    @("77×14=", "89×40="),
    @("24×98=", "68×78="),
    @("83×75=", "91×60="),
    @("70×69=", "69×39="),
    @("18×87=", "96×65="),
    @("62×34=", "27×75="),
    @("65×51=", "77×48=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
